$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commercial")
$ws.Activate()

# Update the two test-data values (Application No / Consumer Number)
# on row 2, columns AN/AO.
$ws.Range("AN2").Value = "JP30000195"
$ws.Range("AO2").Value = "JP30000097"

# Reselect AN:AO (matches the saved selection state in the edited file).
$ws.Range("AN:AO").Select()
